$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "lamenting"
$ws.Range("C3").Value = "To express grief; to weep or wail; to mourn.;To feel great sorrow or regret; to bewail."
$ws.Range("D3").Value = "оплакивать"
